# Refresh the watchlist table with a new batch of NSE tickers and shrink the
# used range from A1:F13 down to A1:F11 (the last two rows are removed).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = "NSE:AAKASH"
$ws.Range("D2").Value = "NSE:BOSCHLTD"
$ws.Range("E2").Value = "NSE:MARUTI"
$ws.Range("F2").Value = "NSE:CUB"

# Row 3
$ws.Range("B3").Value = "NSE:LODHA"
$ws.Range("C3").Value = "NSE:AMBIKCO"
$ws.Range("D3").Value = "NSE:CUB"
$ws.Range("E3").Value = ""

# Row 4
$ws.Range("B4").Value = ""
$ws.Range("C4").Value = "NSE:ATGL"
$ws.Range("D4").Value = "NSE:MUTHOOTFIN"
$ws.Range("E4").Value = ""

# Row 5
$ws.Range("C5").Value = "NSE:CHENNPETRO"
$ws.Range("D5").Value = "NSE:NATIONALUM"
$ws.Range("E5").Value = ""

# Row 6
$ws.Range("C6").Value = "NSE:ELGIEQUIP"
$ws.Range("D6").Value = "NSE:RAMCOCEM"
$ws.Range("E6").Value = ""

# Row 7
$ws.Range("C7").Value = "NSE:GUJALKALI"
$ws.Range("E7").Value = ""

# Row 8
$ws.Range("C8").Value = "NSE:KHADIM"
$ws.Range("E8").Value = ""

# Row 9
$ws.Range("C9").Value = "NSE:LTIM"
$ws.Range("E9").Value = ""

# Row 10
$ws.Range("C10").Value = "NSE:MOL"
$ws.Range("E10").Value = ""

# Row 11
$ws.Range("C11").Value = "NSE:REPRO"
$ws.Range("E11").Value = ""

# Rows 12 and 13 are removed entirely (used range shrinks to A1:F11).
$ws.Range("A12:F13").Delete()
